# Apply predictions.xlsx update:
# - insert 3 new "Colon Adenocarcinoma" rows at the top of the data
# - drop the 2 "Lung Benign Tissue" rows
# - update the "test" folder path to "images for test" for every row
# - keep the remaining Lung rows in their original relative order

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("colonca15.jpeg",  "Colon Adenocarcinoma",          "100.0%", "colonca15.jpeg"),
    @("colonca22.jpeg",  "Colon Adenocarcinoma",          "94.08%", "colonca22.jpeg"),
    @("colonca39.jpeg",  "Colon Adenocarcinoma",          "94.24%", "colonca39.jpeg"),
    @("lungaca122.jpeg", "Lung Adenocarcinoma",           "98.87%", "lungaca122.jpeg"),
    @("lungaca174.jpeg", "Lung Adenocarcinoma",           "99.99%", "lungaca174.jpeg"),
    @("lungaca275.jpeg", "Lung Adenocarcinoma",           "98.69%", "lungaca275.jpeg"),
    @("lungscc315.jpeg", "Lung Squamous Cell Carcinoma",  "100.0%", "lungscc315.jpeg"),
    @("lungscc317.jpeg", "Lung Squamous Cell Carcinoma",  "100.0%", "lungscc317.jpeg")
)

$basePath = "L:/!school/!uni/!classes/sem2-2023/software technology/assignments/assignment 2/images for test/"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $basePath + $data[$i][3]
}
